# Weekly fruit/vegetable price update: insert two new daily records into the
# "Vega Monumental Concepción - Coliflor" dataset.
#
# The dataset is a flat table (row 1 = headers, rows 2..267 = records) with
# constant columns A,B,C,E,F,G,H,N,Q,R for every record in this sheet:
#   A=Mercado ID, B=Mercado, C=Región, E=Codreg, F=Categoría ID,
#   G=Categoría, H=Variedad, N=Unidad de comercialización, Q=Kg o Unidades,
#   R=Clasificación
# and variable columns D=Fecha, I=Calidad, J=Volumen, K=Precio mínimo,
# L=Precio máximo, M=Precio promedio ponderado, O=Origen, P=Precio $/Kg.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Insert #1: new record becomes row 240, pushing the old rows 240.. down ----
$ws.Rows.Item(240).Insert()

$ws.Cells.Item(240, 1).Value = 11
$ws.Cells.Item(240, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(240, 3).Value = "Bíobío"
$ws.Cells.Item(240, 4).Value = 44748
$ws.Cells.Item(240, 5).Value = 8
$ws.Cells.Item(240, 6).Value = 100112008
$ws.Cells.Item(240, 7).Value = "Coliflor"
$ws.Cells.Item(240, 8).Value = "Sin especificar"
$ws.Cells.Item(240, 9).Value = "Primera"
$ws.Cells.Item(240, 10).Value = 500
$ws.Cells.Item(240, 11).Value = 750
$ws.Cells.Item(240, 12).Value = 800
$ws.Cells.Item(240, 13).Value = 780
$ws.Cells.Item(240, 14).Value = "`$/unidad"
$ws.Cells.Item(240, 15).Value = "Región Metropolitana"
$ws.Cells.Item(240, 16).Value = 780
$ws.Cells.Item(240, 17).Value = 1
$ws.Cells.Item(240, 18).Value = "Hortaliza"

# ---- Insert #2: new record becomes row 259, pushing the rest down by one more ----
$ws.Rows.Item(259).Insert()

$ws.Cells.Item(259, 1).Value = 11
$ws.Cells.Item(259, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(259, 3).Value = "Bíobío"
$ws.Cells.Item(259, 4).Value = 44747
$ws.Cells.Item(259, 5).Value = 8
$ws.Cells.Item(259, 6).Value = 100112008
$ws.Cells.Item(259, 7).Value = "Coliflor"
$ws.Cells.Item(259, 8).Value = "Sin especificar"
$ws.Cells.Item(259, 9).Value = "Primera"
$ws.Cells.Item(259, 10).Value = 2500
$ws.Cells.Item(259, 11).Value = 900
$ws.Cells.Item(259, 12).Value = 1000
$ws.Cells.Item(259, 13).Value = 960
$ws.Cells.Item(259, 14).Value = "`$/unidad"
$ws.Cells.Item(259, 15).Value = "Región Metropolitana"
$ws.Cells.Item(259, 16).Value = 960
$ws.Cells.Item(259, 17).Value = 1
$ws.Cells.Item(259, 18).Value = "Hortaliza"

# Make sure the date cells use the same numeric date format as the rest of
# column D (style carries over from Insert(), but set it explicitly too).
$ws.Cells.Item(240, 4).NumberFormat = $ws.Cells.Item(241, 4).NumberFormat
$ws.Cells.Item(259, 4).NumberFormat = $ws.Cells.Item(260, 4).NumberFormat
